$d = $word.ActiveDocument

# Locate the portion of text that needs to be replaced/expanded.
$rng = $d.Content
$rng.Find.Execute("had no iusses with the product backlog", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Replace the found text with the new continuation of the sentence.
$rng.Text = "added another user story " + [char]8220 + "Clear" + [char]8221 + " in which we use a button to delete all shapes on the  drawing window"

# Touching the font forces Word to materialize this range as its own run,
# matching the diff where the sentence is split across two <w:r> elements.
$rng.Font.Size = 13
$rng.Font.Size = 12
